# [LC-943] Update documentation for Letsco OS 1.3.1
# Rename sheets GP1/GP2/BP1..BP7 to GP01/GP02/BP01..BP07 and update the
# matching "KPI ..." title text in cell A1 of each sheet.

$wb = $excel.ActiveWorkbook

$renames = @(
    @{ Old = "GP1"; New = "GP01"; OldTitle = "KPI GP1 - Global Perf 1";        NewTitle = "KPI GP01 - Global Perf 1" },
    @{ Old = "GP2"; New = "GP02"; OldTitle = "KPI GP2 - Global Perf 2";        NewTitle = "KPI GP02 - Global Perf 2" },
    @{ Old = "BP1"; New = "BP01"; OldTitle = "KPI BP1 - Business Process 1";   NewTitle = "KPI BP01 - Business Process 1" },
    @{ Old = "BP2"; New = "BP02"; OldTitle = "KPI BP2 - Business Process 2";   NewTitle = "KPI BP02 - Business Process 2" },
    @{ Old = "BP3"; New = "BP03"; OldTitle = "KPI BP3 - Business Process 3";   NewTitle = "KPI BP03 - Business Process 3" },
    @{ Old = "BP4"; New = "BP04"; OldTitle = "KPI BP4 - Business Process 4";   NewTitle = "KPI BP04 - Business Process 4" },
    @{ Old = "BP5"; New = "BP05"; OldTitle = "KPI BP5 - Business Process 5";   NewTitle = "KPI BP05 - Business Process 5" },
    @{ Old = "BP6"; New = "BP06"; OldTitle = "KPI BP6 - Business Process 6";   NewTitle = "KPI BP06 - Business Process 6" },
    @{ Old = "BP7"; New = "BP07"; OldTitle = "KPI BP7 - Business Process 7";   NewTitle = "KPI BP07 - Business Process 7" }
)

foreach ($entry in $renames) {
    $sheet = $wb.Worksheets.Item($entry.Old)
    $sheet.Range("A1").Value = $entry.NewTitle
    $sheet.Name = $entry.New
}
